$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Step 1: Testable column (B) for rows 2-111: 'y' -> 'n' ---
$ws.Range("B2:B111").Value = 'n'

# --- Step 2: duplicate rows 112:113 down to 114:115 (carries styles + values) ---
$ws.Range("A112:N113").Copy($ws.Range("A114:N115"))
# remove phantom empty H/I/L cells created by the rectangular paste
$ws.Range("H114:I114").Clear()
$ws.Range("L114").Clear()
$ws.Range("H115:I115").Clear()
$ws.Range("L115").Clear()

# --- Step 3: rewrite rows 112 & 113 with the new empty-table-query test cases ---
$ws.Range("G112").Clear()
$ws.Range("B112").Value = 'y'
$ws.Range("C112").Value = '复合索引表空表查询前置'
$ws.Range("F112").Value = 'mixindex010'
$ws.Range("I112").Value = 'select /*+ vector_pre */ id,name,age,feature_id,feature_index$distance from vector($mixindex010, feature, array[1.3893274068832397, 0.4223838150501251, 0.7195155024528503, 0.28465819358825684, 0.9267012476921082, 0.4375186264514923, 0.7362583875656128, 0.41960853338241577, 0.7168405055999756, 0.5768887996673584, 0.4182721674442291, 0.129996195435524, 0.21704305708408356, 0.35628095269203186, 0.5203919410705566, 0.4046420454978943, 0.1657610833644867, 0.7787348031997681, 0.16017264127731323, 0.7188393473625183, 0.916609525680542, 0.6900423765182495, 0.6827380657196045, 0.49170464277267456, 0.5505375266075134, 0.33907604217529297, 0.09597073495388031, 0.008679530583322048, 0.9524646997451782, 0.2128734141588211, 0.6533687710762024, 0.6094813346862793, 0.18400055170059204, 0.9364618062973022, 0.7497748732566833, 0.05966084823012352, 0.40710607171058655, 0.6920192241668701, 0.649844765663147, 0.705480694770813, 0.11741353571414948, 0.2164693921804428, 0.48223137855529785, 0.6425648331642151, 0.5236963033676147, 0.6490180492401123, 0.5132198333740234, 0.5319958925247192, 0.5141375660896301, 0.9233165979385376, 0.5702359080314636, 0.46916520595550537, 0.7964460849761963, 0.17150050401687622, 0.4680892825126648, 0.6861740946769714, 0.4802965819835663, 0.9073042273521423, 0.17548426985740662, 0.9943628311157227, 0.9149019718170166, 0.7661579847335815, 0.13299474120140076, 0.7843778133392334], 10, map[efSearch, 40]) where age=78 and name=''Wd'' order by feature_index$distance limit 10'
$ws.Range("K112").Value = 'csv_equals'
$ws.Range("L112").Value = 'explain plan for select /*+ vector_pre */ id,name,age,feature_id,feature_index$distance from vector($mixindex010, feature, array[1.3893274068832397, 0.4223838150501251, 0.7195155024528503, 0.28465819358825684, 0.9267012476921082, 0.4375186264514923, 0.7362583875656128, 0.41960853338241577, 0.7168405055999756, 0.5768887996673584, 0.4182721674442291, 0.129996195435524, 0.21704305708408356, 0.35628095269203186, 0.5203919410705566, 0.4046420454978943, 0.1657610833644867, 0.7787348031997681, 0.16017264127731323, 0.7188393473625183, 0.916609525680542, 0.6900423765182495, 0.6827380657196045, 0.49170464277267456, 0.5505375266075134, 0.33907604217529297, 0.09597073495388031, 0.008679530583322048, 0.9524646997451782, 0.2128734141588211, 0.6533687710762024, 0.6094813346862793, 0.18400055170059204, 0.9364618062973022, 0.7497748732566833, 0.05966084823012352, 0.40710607171058655, 0.6920192241668701, 0.649844765663147, 0.705480694770813, 0.11741353571414948, 0.2164693921804428, 0.48223137855529785, 0.6425648331642151, 0.5236963033676147, 0.6490180492401123, 0.5132198333740234, 0.5319958925247192, 0.5141375660896301, 0.9233165979385376, 0.5702359080314636, 0.46916520595550537, 0.7964460849761963, 0.17150050401687622, 0.4680892825126648, 0.6861740946769714, 0.4802965819835663, 0.9073042273521423, 0.17548426985740662, 0.9943628311157227, 0.9149019718170166, 0.7661579847335815, 0.13299474120140076, 0.7843778133392334], 10, map[efSearch, 40]) where age=78 and name=''Wd'' order by feature_index$distance limit 10'
$ws.Range("G113").Clear()
$ws.Range("B113").Value = 'y'
$ws.Range("C113").Value = '复合索引表空表查询后置'
$ws.Range("F113").Value = 'mixindex010'
$ws.Range("I113").Value = 'select name,age,feature_id,feature_index$distance from vector($mixindex010, feature, array[1.3893274068832397, 0.4223838150501251, 0.7195155024528503, 0.28465819358825684, 0.9267012476921082, 0.4375186264514923, 0.7362583875656128, 0.41960853338241577, 0.7168405055999756, 0.5768887996673584, 0.4182721674442291, 0.129996195435524, 0.21704305708408356, 0.35628095269203186, 0.5203919410705566, 0.4046420454978943, 0.1657610833644867, 0.7787348031997681, 0.16017264127731323, 0.7188393473625183, 0.916609525680542, 0.6900423765182495, 0.6827380657196045, 0.49170464277267456, 0.5505375266075134, 0.33907604217529297, 0.09597073495388031, 0.008679530583322048, 0.9524646997451782, 0.2128734141588211, 0.6533687710762024, 0.6094813346862793, 0.18400055170059204, 0.9364618062973022, 0.7497748732566833, 0.05966084823012352, 0.40710607171058655, 0.6920192241668701, 0.649844765663147, 0.705480694770813, 0.11741353571414948, 0.2164693921804428, 0.48223137855529785, 0.6425648331642151, 0.5236963033676147, 0.6490180492401123, 0.5132198333740234, 0.5319958925247192, 0.5141375660896301, 0.9233165979385376, 0.5702359080314636, 0.46916520595550537, 0.7964460849761963, 0.17150050401687622, 0.4680892825126648, 0.6861740946769714, 0.4802965819835663, 0.9073042273521423, 0.17548426985740662, 0.9943628311157227, 0.9149019718170166, 0.7661579847335815, 0.13299474120140076, 0.7843778133392334], 10, map[efSearch, 40]) where name=''71vF'' or age=22 order by feature_index$distance limit 10'
$ws.Range("K113").Value = 'csv_equals'
$ws.Range("L113").Value = 'explain plan for select name,age,feature_id,feature_index$distance from vector($mixindex010, feature, array[1.3893274068832397, 0.4223838150501251, 0.7195155024528503, 0.28465819358825684, 0.9267012476921082, 0.4375186264514923, 0.7362583875656128, 0.41960853338241577, 0.7168405055999756, 0.5768887996673584, 0.4182721674442291, 0.129996195435524, 0.21704305708408356, 0.35628095269203186, 0.5203919410705566, 0.4046420454978943, 0.1657610833644867, 0.7787348031997681, 0.16017264127731323, 0.7188393473625183, 0.916609525680542, 0.6900423765182495, 0.6827380657196045, 0.49170464277267456, 0.5505375266075134, 0.33907604217529297, 0.09597073495388031, 0.008679530583322048, 0.9524646997451782, 0.2128734141588211, 0.6533687710762024, 0.6094813346862793, 0.18400055170059204, 0.9364618062973022, 0.7497748732566833, 0.05966084823012352, 0.40710607171058655, 0.6920192241668701, 0.649844765663147, 0.705480694770813, 0.11741353571414948, 0.2164693921804428, 0.48223137855529785, 0.6425648331642151, 0.5236963033676147, 0.6490180492401123, 0.5132198333740234, 0.5319958925247192, 0.5141375660896301, 0.9233165979385376, 0.5702359080314636, 0.46916520595550537, 0.7964460849761963, 0.17150050401687622, 0.4680892825126648, 0.6861740946769714, 0.4802965819835663, 0.9073042273521423, 0.17548426985740662, 0.9943628311157227, 0.9149019718170166, 0.7661579847335815, 0.13299474120140076, 0.7843778133392334], 10, map[efSearch, 40]) where name=''71vF'' or age=22 order by feature_index$distance limit 10'

# --- Step 4: set the new-id columns (A) + remaining new-string columns (J, M) on rows 114 & 115 ---
$ws.Range("A114").Value = 'mixindex_113'
$ws.Range("J114").Value = 'src/test/resources/io.dingodb.test/testdata/cases/index/expectedresult/mix/mixindex_113.csv'
$ws.Range("M114").Value = 'src/test/resources/io.dingodb.test/testdata/cases/index/expectedresult/mix/mixindex_113_explain.csv'
$ws.Range("A115").Value = 'mixindex_114'
$ws.Range("J115").Value = 'src/test/resources/io.dingodb.test/testdata/cases/index/expectedresult/mix/mixindex_114.csv'
$ws.Range("M115").Value = 'src/test/resources/io.dingodb.test/testdata/cases/index/expectedresult/mix/mixindex_114_explain.csv'

# --- Step 5: sheet view (selection / scroll position) ---
$ws.Range("H108").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 85
$win.ScrollColumn = 1
